$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at the top; shifts the two existing tables
# (rows 1-5 and rows 17-21) down by one row each.
$ws.Rows("1").Insert()

# The former "Wire Color" column (now sitting in column B, rows 2-6)
# becomes the new "Wire Color between Mirrors" column C -- copy it
# over before column B is overwritten with the new data.
$ws.Range("C2:C6").Value2 = $ws.Range("B2:B6").Value2

# New header row describing the two wire-color columns + box side.
$ws.Range("B1").Value = "test setup"
$ws.Range("C1").Value = "Wire Color between Mirrors"
$ws.Range("D1").Value = "Side of Box"
$ws.Range("E1").Value = "Wire Color in Box"

# New "Wire Color in Box" values for column B.
$ws.Range("B2").Value = "Blue"
$ws.Range("B3").Value = "Purple"
$ws.Range("B4").Value = "Grey"
$ws.Range("B5").Value = "White"
$ws.Range("B6").Value = "Black"

# Widen the two new text columns.
$ws.Columns("C").ColumnWidth = 23.666666666666668
$ws.Columns("E").ColumnWidth = 14.166666666666666

# Restore the selection used in the authored workbook.
$ws.Range("B7").Select() | Out-Null
